# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-50, columns E/F/G) is refreshed:
#   - Column E (period label, e.g. "1705") is rewritten so the period list
#     now runs oldest -> newest (1705 .. 2003) instead of newest -> oldest.
#   - Column F ("Valor Mora") and column G ("Salario Basico") are updated to
#     the new figures that came with this refresh of the account-statement
#     database.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i

    # New "Valor Mora" figures: rows 16-34 take 29509, rows 35-50 take 31249
    # (this mirrors the original split, with the two figures swapped).
    if ($row -le 34) {
        $valorMora = 29509
    } else {
        $valorMora = 31249
    }

    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valorMora
    $ws.Range("G$row").Value = 781242
}
